$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65; this shifts existing rows 65-98 down to 66-99
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new record's data
$ws.Cells.Item(65, 1).Value2 = 4
$ws.Cells.Item(65, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(65, 3).Value2 = "Los Lagos"
$ws.Cells.Item(65, 4).Value2 = 44813
$ws.Cells.Item(65, 5).Value2 = 10
$ws.Cells.Item(65, 6).Value2 = 100112031
$ws.Cells.Item(65, 7).Value2 = "Poroto verde"
$ws.Cells.Item(65, 8).Value2 = "Magnum"
$ws.Cells.Item(65, 9).Value2 = "Primera"
$ws.Cells.Item(65, 10).Value2 = 35
$ws.Cells.Item(65, 11).Value2 = 37000
$ws.Cells.Item(65, 12).Value2 = 37000
$ws.Cells.Item(65, 13).Value2 = 37000
$ws.Cells.Item(65, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(65, 15).Value2 = "Perú"
$ws.Cells.Item(65, 16).Value2 = 1480
$ws.Cells.Item(65, 17).Value2 = 25
$ws.Cells.Item(65, 18).Value2 = "Hortaliza"
